# Apply "IVG" update: drop the trailing "<longueur> <debut> <fin>" position
# hints from the libelle column (F), move the active-cell selection, and
# nudge the sheet's default column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing " <longueur> <debut> <fin>" numeric suffix from the
# libelle (column F) cells.
$ws.Range("F3").Value  = "N° FINESS du fichier d'entrée"
$ws.Range("F4").Value  = "Type de prestation"
$ws.Range("F5").Value  = "Année période"
$ws.Range("F6").Value  = "N° période (mois)"
$ws.Range("F7").Value  = "N° d'index du RSA"
$ws.Range("F8").Value  = "Mois du séjour"
$ws.Range("F9").Value  = "Année du séjour"
$ws.Range("F10").Value = "Nombre d'IVG antérieures"
$ws.Range("F11").Value = "Année de la dernière IVG"
$ws.Range("F12").Value = "Nombre de naissances vivantes antérieures"
$ws.Range("F13").Value = "Filler"

# Move the current selection from H4 to F14.
$ws.Range("F14").Select()

# Nudge the sheet's default column width slightly.
$ws.StandardWidth = 10.5859375
